# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# New K values were (re)computed from the underlying box-score data and are
# written here as literal values, one per game row (rows 2-51 on Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 3
    7  = 0
    8  = 1
    9  = 5
    10 = 2
    11 = 3
    12 = 1
    13 = 1
    14 = 2
    15 = 0
    16 = 2
    17 = 2
    18 = 1
    19 = 3
    20 = 0
    21 = 1
    22 = 3
    23 = 0
    24 = 4
    25 = 3
    26 = 2
    27 = 1
    28 = 2
    29 = 1
    30 = 2
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 0
    38 = 1
    39 = 2
    40 = 1
    41 = 3
    42 = 0
    43 = 0
    44 = 3
    45 = 4
    46 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 2
    51 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
